$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Rearranged headers to be similar to PrePrints": the whole contents of
# column B ("BuSel G") and column G ("Ispanya") are swapped - header label,
# the 17 data rows, and the footer count label in row 21 all move together.
# Row 20 (the "TOPLAM" SUM-formula row) is left alone: SUM(B2:B18) and
# SUM(G2:G18) keep pointing at their own column, so only their cached
# results change once the underlying data moves.
#
# A scratch column (Z) well outside the used range stages column B's
# original contents while column G's contents are copied into B; the
# staged values then move into G. Every destination is cleared right
# before it receives a copy because Range.Copy() is a no-op on a truly
# empty source cell (it would otherwise leave stale contents/format
# behind in the destination).

$valueBlocks = @("B1:B18", "B21")
$scratchBlocks = @("Z1:Z18", "Z21")

for ($i = 0; $i -lt $valueBlocks.Length; $i++) {
    $ws.Range($scratchBlocks[$i]).Clear()
    $ws.Range($valueBlocks[$i]).Copy($ws.Range($scratchBlocks[$i]))
}

$otherBlocks = @("G1:G18", "G21")

for ($i = 0; $i -lt $valueBlocks.Length; $i++) {
    $ws.Range($valueBlocks[$i]).Clear()
    $ws.Range($otherBlocks[$i]).Copy($ws.Range($valueBlocks[$i]))
}

for ($i = 0; $i -lt $valueBlocks.Length; $i++) {
    $ws.Range($otherBlocks[$i]).Clear()
    $ws.Range($scratchBlocks[$i]).Copy($ws.Range($otherBlocks[$i]))
}

$ws.Range("Z1:Z21").Clear()

# The new column-B header (formerly "Ispanya", now holding what used to be
# column G's data) is relabelled "Andalousia".
$ws.Range("B1").Value = "Andalousia"
